# Apply the "BAL added" update to the Bill_data_summary workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bill rows (9-14) appended beneath the existing data table.
# Values are entered column-by-column (as the original author did),
# which also determines the order new entries land in the shared
# string table.
$billIds = @("BILL008", "BILL009", "BILL010", "BILL011", "BILL012", "BILL013")
$billFor = @("Richard", "Gibson", "Richard", "Nagalingaraj", "<NULL>", "Richard")
$amounts = @(2344.67, 144.56, 2000, 1500, $null, 7544.08)
$paidBy  = @("Card", "Card", "Cash", "Cash", $null, "Card")

for ($i = 0; $i -lt $billIds.Length; $i++) {
    $ws.Cells.Item(9 + $i, 1).Value = $billIds[$i]
}
for ($i = 0; $i -lt $billFor.Length; $i++) {
    $ws.Cells.Item(9 + $i, 2).Value = $billFor[$i]
}
for ($i = 0; $i -lt $paidBy.Length; $i++) {
    if ($paidBy[$i] -ne $null) {
        $ws.Cells.Item(9 + $i, 4).Value = $paidBy[$i]
    }
}
for ($i = 0; $i -lt $amounts.Length; $i++) {
    if ($amounts[$i] -ne $null) {
        $ws.Cells.Item(9 + $i, 3).Value = $amounts[$i]
    }
}

# Row 13's amount is the corrupt/null placeholder text rather than a number.
$ws.Cells.Item(13, 3).Value = "<Bnull>"

# Column B was autosized ("best fit") by Excel after the paste.
$ws.Columns.Item(2).ColumnWidth = 11.109375

# Selection left on F8 after the refresh (per the author's note, the
# refresh bug was not fixed).
$ws.Range("F8").Select()
